# The sheet gained one new data row: a new record was inserted as row 354
# (pushing the former rows 354-426 down to 355-427). Replicate that with a
# real row insert so every subsequent row's data shifts down automatically,
# then populate the freshly inserted row with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 354; Excel shifts rows 354:426 down to 355:427.
$ws.Range("A354").EntireRow.Insert()

# Populate the newly inserted row 354 with the new record's data.
$ws.Range("A354").Value = 9
$ws.Range("B354").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C354").Value = "Metropolitana"
$ws.Range("D354").Value = 44711
$ws.Range("E354").Value = 13
$ws.Range("F354").Value = "Fruta"
$ws.Range("G354").Value = 100108
$ws.Range("H354").Value = "Tropicales y subtropicales"
$ws.Range("I354").Value = 100108002
$ws.Range("J354").Value = "Mango"
$ws.Range("K354").Value = "Sin especificar"
$ws.Range("L354").Value = "Especial"
$ws.Range("M354").Value = 120
$ws.Range("N354").Value = 10000
$ws.Range("O354").Value = 10000
$ws.Range("P354").Value = 10000
$ws.Range("Q354").Value = "`$/bandeja 4 kilos"
$ws.Range("R354").Value = "Brasil"
$ws.Range("S354").Value = 2500
$ws.Range("T354").Value = 4
